{"js": "const oldText = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie ozvezdje Lev 2022: 14.-23. april, 14.-23. maj\";\nconst newText = \"2022: Datumi kampanje za opazovanje ozvezdje Lev: 14.-23. april, 14.-23. maj\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie ozvezdje Lev 2022: 14.-23. april, 14.-23. maj\"\n$newText = \"2022: Datumi kampanje za opazovanje ozvezdje Lev: 14.-23. april, 14.-23. maj\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n"}
